# "included total stack frames we're using for training"
# Adds two new columns to the "After keeping max 26 vids" sheet:
#   U: "Num Words" - number of words in each sentence (column B)
#   W: "Training Size" - U * (sum of the P1..P6 "# Training" counts, C:H)
# plus a "Total" row at the bottom of the new W column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- Header row (row 1) ----
$ws.Range("U1").Value = "Num Words"
$ws.Range("U1").Font.Bold = $true

$ws.Range("W1").Value = "Training Size"
$ws.Range("W1").Font.Bold = $true

# ---- Per-sentence data (rows 2-51) ----
# Number of words per sentence (column B), entered as the literal word count.
$numWords = @{
    2=2; 3=2; 4=4; 5=5; 6=5; 7=5; 8=5; 9=3; 10=3; 11=2;
    12=4; 13=3; 14=2; 15=2; 16=3; 17=3; 18=2; 19=4; 20=3; 21=3;
    22=4; 23=2; 24=3; 25=5; 26=5; 27=5; 28=2; 29=2; 30=4; 31=4;
    32=5; 33=4; 34=3; 35=2; 36=3; 37=3; 38=3; 39=3; 40=2; 41=2;
    42=2; 43=3; 44=2; 45=2; 46=2; 47=2; 48=4; 49=2; 50=2; 51=3
}

for ($row = 2; $row -le 51; $row++) {
    $ws.Range("U$row").Value = $numWords[$row]
    $ws.Range("U$row").HorizontalAlignment = -4108   # xlCenter

    $ws.Range("W$row").Formula = "=U$row*(C$row+D$row+E$row+F$row+G$row+H$row)"
}

# ---- Totals row (row 52) ----
$ws.Range("I52").Copy() | Out-Null
$ws.Range("V52").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$ws.Range("V52").Value = "Total"

$ws.Range("J52").Copy() | Out-Null
$ws.Range("W52").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$ws.Range("W52").Formula = "=SUM(W2:W51)"

# ---- Column width for the new "Training Size" column ----
$ws.Columns.Item(23).ColumnWidth = 10.7

# ---- View / selection (best effort - scroll position) ----
$ws.Activate()
$ws.Range("S53").Select()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 7

# ---- Page setup ----
$ws.PageSetup.Orientation = 1   # xlPortrait
